$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - update "Use" label only (keep as text, preserve leading space)
$ws.Range("A2").Value = "'" + " 46"
$ws.Range("A2").Style = "Normal"

# Row 3 - update all stats
$ws.Range("A3").Value = "'" + " 31"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 91
$ws.Range("C3").Value = 0.84
$ws.Range("D3").Value = 157.5
$ws.Range("E3").Value = 0.73
$ws.Range("F3").Value = 98.90000000000001
$ws.Range("G3").Value = 0.6899999999999999
$ws.Range("H3").Value = 0.54
$ws.Range("I3").Value = 0.37
$ws.Range("J3").Value = 0.03
$ws.Range("K3").Value = 0.08
$ws.Range("L3").Value = 49
$ws.Range("M3").Value = 67
$ws.Range("N3").Value = 34
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 7

# Row 4 - update "Use" label only (keep as text, preserve leading space)
$ws.Range("A4").Value = "'" + " 15"
$ws.Range("A4").Style = "Normal"
